# Add a new "time_taken" column (F) to the sheet, with a header in F1
# (styled like the other header cells B1:E1) and per-row timestamp
# values in F2:F7, matching the data rows' (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("F1").Value = "time_taken"

# Copy the header style (bold font, border, centered/top alignment)
# from an existing header cell so F1 matches B1:E1.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row timestamps (plain text, same as the other data columns)
$ws.Range("F2").Value = "2021-10-05 13:38:38.438841"
$ws.Range("F3").Value = "2021-10-05 13:38:38.438852"
$ws.Range("F4").Value = "2021-10-05 13:38:38.438855"
$ws.Range("F5").Value = "2021-10-05 13:38:38.438858"
$ws.Range("F6").Value = "2021-10-05 13:38:38.438861"
$ws.Range("F7").Value = "2021-10-05 13:38:38.438864"
